$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the values in column A (IP_Hoa Binh -> IP_BRVT)
$ws.Range("A2").Value = "IP_BRVT"
$ws.Range("A3").Value = "IP_BRVT"

# Column B content is no longer present for rows 2-3
$ws.Range("B2").ClearContents()
$ws.Range("B3").ClearContents()

# Column C now carries District / Sub_Type instead of Developer
$ws.Range("C2").Value = "District"
$ws.Range("C3").Value = "Sub_Type"

# Row 4 is removed entirely
$ws.Rows.Item(4).Delete()
